$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D-column values are stored as plain text (some contain multiple
# "." separators, e.g. "29.292.50", which Excel would otherwise try
# to coerce into a number). Forcing a text number format before the
# assignment keeps the literal string, then resetting the style back
# to "Normal" avoids leaving a stray style applied to the cell.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "29.292.50"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.43%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.918.82"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  -0.19%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "328.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("E6").Value = "  -0.22%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4666"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3964"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.56%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "47.11"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.08041"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.020"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.67%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "22.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.915.79"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.159"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("E15").Value = "  +0.81%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.06976"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "89.97"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.00001019"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.17%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "29.280.29"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.389"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  +0.49%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.136.84"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.066"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.56%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "155.54"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.81"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.917"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.16%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.037"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.42%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "121.24"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("E32").Value = "  -0.13%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9455"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.35%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.381"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.367"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.263"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05881"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.169"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.24%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.179"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.02112"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5864"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "

$ws.Range("E42").Value = "  -0.18%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1823"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("E44").Value = "  +1.18%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.322"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +8.74%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5488"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "12.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.07237"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.895"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.138"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "113.51"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "

